$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.301.99'
$ws.Range('E2').Value = '  +7.04%  '
$ws.Range('D3').Value = '3.686.44'
$ws.Range('E3').Value = '  +19.33%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''598.27'
$ws.Range('E5').Value = '  +3.90%  '
$ws.Range('D6').Value = '''184.25'
$ws.Range('E6').Value = '  +6.94%  '
$ws.Range('D7').Value = '3.686.05'
$ws.Range('E7').Value = '  +19.38%  '
$ws.Range('D8').Value = '''0.999'
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').Value = '''0.536'
$ws.Range('E9').Value = '  +4.40%  '
$ws.Range('D10').Value = '''0.164'
$ws.Range('E10').Value = '  +8.43%  '
$ws.Range('D11').Value = '''6.62'
$ws.Range('E11').Value = '  +4.21%  '
$ws.Range('E12').Value = '  +6.48%  '
$ws.Range('D13').Value = '''39.93'
$ws.Range('E13').Value = '  +11.94%  '
$ws.Range('D14').Value = '''0.0000255'
$ws.Range('E14').Value = '  +6.76%  '
$ws.Range('D15').Value = '4.301.59'
$ws.Range('E15').Value = '  +19.25%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.684.61'
$ws.Range('E16').Value = '  +19.15%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '71.283.96'
$ws.Range('E17').Value = '  +7.05%  '
$ws.Range('E18').Value = '  +1.83%  '
$ws.Range('D19').Value = '''7.50'
$ws.Range('E19').Value = '  +7.86%  '
$ws.Range('D20').Value = '''16.93'
$ws.Range('E20').Value = '  +1.11%  '
$ws.Range('D21').Value = '''516.17'
$ws.Range('E21').Value = '  +6.82%  '
$ws.Range('D22').Value = '''9.20'
$ws.Range('E22').Value = '  +18.49%  '
$ws.Range('D23').Value = '''0.746'
$ws.Range('E23').Value = '  +8.78%  '
$ws.Range('D24').Value = '''87.65'
$ws.Range('E24').Value = '  +5.32%  '
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').Value = '''2.42'
$ws.Range('E25').Value = '  +8.64%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Value = '''13.47'
$ws.Range('E26').Value = '  +6.55%  '
$ws.Range('D27').Value = '''10.88'
$ws.Range('E27').Value = '  +8.28%  '
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('D29').Value = '''2.53'
$ws.Range('E29').Value = '  +12.74%  '
$ws.Range('D30').Value = '''8.18'
$ws.Range('E30').Value = '  +3.16%  '
$ws.Range('D31').Value = '''31.78'
$ws.Range('E31').Value = '  +13.96%  '
$ws.Range('D32').Value = '''2.77'
$ws.Range('E32').Value = '  +7.18%  '
$ws.Range('E33').Value = '  +17.46%  '
$ws.Range('E34').Value = '  +4.65%  '
$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('D36').Value = '''6.15'
$ws.Range('E36').Value = '  +10.57%  '
$ws.Range('E37').Value = '  +8.15%  '
$ws.Range('E38').Value = '  +10.59%  '
$ws.Range('D39').Value = '''2.15'
$ws.Range('E39').Value = '  +9.90%  '
$ws.Range('D40').Value = '''50.82'
$ws.Range('E40').Value = '  +3.71%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '''0.128'
$ws.Range('E41').Value = '  +4.36%  '
$ws.Range('B42').Value = 'Arweave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D42').Value = '''45.78'
$ws.Range('E42').Value = '  -5.78%  '
$ws.Range('D43').Value = '3.177.17'
$ws.Range('E43').Value = '  +14.38%  '
$ws.Range('D44').Value = '''8.81'
$ws.Range('E44').Value = '  +6.92%  '
$ws.Range('D45').Value = '''2.78'
$ws.Range('E45').Value = '  +7.09%  '
$ws.Range('D46').Value = '''406.76'
$ws.Range('E46').Value = '  +10.59%  '
$ws.Range('D47').Value = '''0.0368'
$ws.Range('E47').Value = '  +6.55%  '
$ws.Range('D48').Value = '''28.13'
$ws.Range('E48').Value = '  +15.39%  '
$ws.Range('D49').Value = '''136.56'
$ws.Range('E49').Value = '  +1.81%  '
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('E51').Value = '  +12.68%  '
